$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "TextBoxInvalidEmail" sheet (3rd sheet): add 3 new rows of test data
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A3").Value = "peracom"

$ws3.Range("A4").Value = "/pera@pera.com"
$null = $ws3.Hyperlinks.Add($ws3.Range("A4"), "mailto:/pera@pera.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "/pera@pera.com")
# Restyle the link cell so it matches the workbook's existing custom
# "hyperlink" look (blue text, no underline) instead of Excel's builtin style.
$ws3.Range("A4").ClearFormats()
$ws3.Range("A4").Value = "/pera@pera.com"
$ws3.Range("A4").Font.Color = 16711680

$ws3.Range("A5").Value = "pera.com@"

$ws3.Columns("A").ColumnWidth = 15.8

# ---------------------------------------------------------------------------
# 2. New sheet "BookStoreLogIn" (added after "TextBoxInvalidEmail")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "BookStoreLogIn"

$ws4.Range("A1").Value = "ValidUsername"
$ws4.Range("B1").Value = "ValidPassword"

$ws4.Range("A2").Value = "helenatodorovic86@gmail.com"
$null = $ws4.Hyperlinks.Add($ws4.Range("A2"), "mailto:helenatodorovic86@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "helenatodorovic86@gmail.com")
$ws4.Range("A2").ClearFormats()
$ws4.Range("A2").Value = "helenatodorovic86@gmail.com"
$ws4.Range("A2").Font.Color = 16711680

$ws4.Range("B2").Value = "ITBootcamp2023!"

$ws4.Columns("A").ColumnWidth = 26.8
$ws4.Columns("B").ColumnWidth = 14.3

$null = $ws4.Range("B2").Select()
$null = $ws4.Activate()
